$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the full previous data block (A2:E38) so no stale values are left behind
# (e.g. row 22 used to have B/C/E populated but the refreshed "430" row does not).
$ws.Range("A2:E38").ClearContents()

# Column A holds the "evento" codes as text (matches source formatting), so force
# text number format before writing the values.
$ws.Range("A2:A35").NumberFormat = "@"

$data = @(
    @("100", "Accidente ofidico", 0, 0, 1),
    @("113", "Desnutrici”n aguda en menores de 5 anos", 1, 2, 0.18),
    @("115", "Cancer en menores de 18 anos", 0, 0, 1),
    @("155", "Cancer de la mama y cuello uterino", 8, 4, 0.06),
    @("210", "Dengue", 4, 27, 0),
    @("215", "Defectos congenitos", 2, 2, 0.27),
    @("217", "Chikungunya", 0, 0, 1),
    @("220", "Dengue grave", 0, 0, 1),
    @("298", "Evento adverso grave posterior a la vacunacion", 0, 0, 1),
    @("300", "Agresiones por animales potencialmente transmisores de rabia", 45, 20, 0),
    @("330", "Hepatitis a", 2, 0, 0.14),
    @("340", "Hepatitis b, c y coinfeccion hepatitis b y delta", 1, 2, 0.18),
    @("342", "Enfermedades huerfanas - raras", 5, 5, 0.18),
    @("348", "Infeccion respiratoria aguda grave irag inusitada", 2, 0, 0.14),
    @("352", "Infecciones de sitio quirurgico asociadas a procedimiento medico quirurgico", 1, 0, 0.37),
    @("355", "Enfermedad transmitida por alimentos o agua (eta)", 0, 0, 1),
    @("356", "Intento de suicidio", 7, 8, 0.13),
    @("357", "Iad - infecciones asociadas a dispositivos - individual", 2, 0, 0.14),
    @("365", "Intoxicaciones", 5, 4, 0.18),
    @("420", "Leishmaniasis cutanea", 0, 0, 1),
    @("430", $null, $null, 0, $null),
    @("455", "Leptospirosis", 1, 1, 0.37),
    @("465", "Malaria", 0, 0, 1),
    @("535", "Meningitis bacteriana y enfermedad meningoc”cica", 0, 0, 1),
    @("549", "Morbilidad materna extrema", 6, 3, 0.09),
    @("560", "Mortalidad perinatal y neonatal tardia", 1, 2, 0.18),
    @("580", "Mortalidad por dengue", 0, 0, 1),
    @("620", "Parotiditis", 1, 0, 0.37),
    @("740", "Sifilis congenita", 1, 0, 0.37),
    @("750", "Sifilis gestacional", 2, 0, 0.14),
    @("813", "Tuberculosis", 7, 2, 0.02),
    @("831", "Varicela individual", 5, 2, 0.08),
    @("850", "Vih/sida/mortalidad por sida", 10, 7, 0.09),
    @("895", "Zika", 0, 0, 1)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $c + 1).Value = $val
        }
    }
}

# Restore default styling on column A now that the text type is recorded --
# NumberFormat="@" above only exists to force text storage, it is not part of the
# original look of the sheet.
$ws.Range("A2:A35").Style = "Normal"

$ws.Range("A1").Select()
